# Update attendance figures ("想去人数", column F) for a handful of events
# on the "展览" and "全部类型" sheets, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value = 7657
        $ws.Range("F5").Value = 39
        $ws.Range("F6").Value = 279
        $ws.Range("F7").Value = 1147
        $ws.Range("F8").Value = 202
        $ws.Range("F10").Value = 156
    }
    else {
        $ws.Range("F2").Value = 7657
        $ws.Range("F5").Value = 39
        $ws.Range("F6").Value = 279
        $ws.Range("F7").Value = 1147
        $ws.Range("F8").Value = 202
        $ws.Range("F11").Value = 156
    }
}

$wb.Save()
